# Golang_SHA256_Benchmark_CPU_Test.xlsx - "Add files via upload"
# Refresh the benchmark table: a few CPU / clock-speed / RAM labels were
# corrected, two brand-new systems (ANTMINER S9 row 18, Raspberry Pi B row 21)
# were added, and the previously-empty row 18 (ANTMINER S9) now carries its
# benchmark time + date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row ----------------------------------------------------------
$ws.Range("A1").Value = "System"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "CPU"
$ws.Range("D1").Value = "GHz"
$ws.Range("E1").Value = "RAM"
$ws.Range("F1").Value = "1M (ms)"
$ws.Range("G1").Value = "Date"

# ---- data rows -------------------------------------------------------------
$ws.Range("A2").Value = "Windows10"
$ws.Range("B2").Value = "PC"
$ws.Range("C2").Value = "Ryzen 5 3600"
$ws.Range("D2").Value = "3.6"
$ws.Range("E2").Value = "DDR4-3200 16GB"
$ws.Range("F2").Value = 85
$ws.Range("G2").Value = 45323

$ws.Range("A3").Value = "Mac"
$ws.Range("B3").Value = "MacBook Pro A1398"
$ws.Range("C3").Value = "Intel Core i7"
$ws.Range("D3").Value = "2.8"
$ws.Range("E3").Value = "DDR4 16GB"
$ws.Range("F3").Value = 250
$ws.Range("G3").Value = 45323

$ws.Range("A4").Value = "Windows10"
$ws.Range("B4").Value = "Notebook TingPad Yoda460"
$ws.Range("C4").Value = "Intel i5-6200U"
$ws.Range("D4").Value = "2.8"
$ws.Range("E4").Value = "DDR4 8GB"
$ws.Range("F4").Value = 420
$ws.Range("G4").Value = 45323

$ws.Range("A5").Value = "Windows10"
$ws.Range("B5").Value = "PC"
$ws.Range("C5").Value = "AMD Phenom II X4 975"
$ws.Range("D5").Value = "3.6"
$ws.Range("E5").Value = "DDR3 16GB"
$ws.Range("F5").Value = 477
$ws.Range("G5").Value = 45324

$ws.Range("A6").Value = "Linux"
$ws.Range("B6").Value = "VPS Server Contabo"
$ws.Range("C6").Value = "Intel i5"
$ws.Range("D6").Value = "2.8"
$ws.Range("E6").Value = "DDR4 8GB"
$ws.Range("F6").Value = 500
$ws.Range("G6").Value = 45323

$ws.Range("A7").Value = "Termux Android "
$ws.Range("B7").Value = "Redmi Note 8 Pro"
$ws.Range("C7").Value = "Helio G90T"
$ws.Range("D7").Value = "2.0"
$ws.Range("E7").Value = "DDR4 6GB"
$ws.Range("F7").Value = 531
$ws.Range("G7").Value = 45323
$ws.Range("H7").Value = "186-877"

$ws.Range("A8").Value = "Windows10"
$ws.Range("B8").Value = "PC Dell"
$ws.Range("C8").Value = "Intel E8500 Duo"
$ws.Range("D8").Value = "3.16"
$ws.Range("E8").Value = "DDR3 8GB"
$ws.Range("F8").Value = 547
$ws.Range("G8").Value = 45323

$ws.Range("A9").Value = "Termux Android "
$ws.Range("B9").Value = "Mi Pad 4 Plus"
$ws.Range("C9").Value = "Snapdragon 660"
$ws.Range("D9").Value = "2.2"
$ws.Range("E9").Value = "DDR4 4GB"
$ws.Range("F9").Value = 627
$ws.Range("G9").Value = 45323
$ws.Range("H9").Value = "355-900"

$ws.Range("A10").Value = "Linux"
$ws.Range("B10").Value = "NanoPi NEO3"
$ws.Range("C10").Value = "RockChip RK3328"
$ws.Range("D10").Value = "1.3"
$ws.Range("E10").Value = "DDR4 2GB"
$ws.Range("F10").Value = 715
$ws.Range("G10").Value = 45323

$ws.Range("A11").Value = "Linux"
$ws.Range("B11").Value = "Orange Pi Zero3"
$ws.Range("C11").Value = "Allwinner H618"
$ws.Range("D11").Value = "1.5"
$ws.Range("E11").Value = "DDR4 4GB"
$ws.Range("F11").Value = 990
$ws.Range("G11").Value = 45323

$ws.Range("A12").Value = "Windows7"
$ws.Range("B12").Value = "Notebook Lenovo G500"
$ws.Range("C12").Value = "Intel Celeron 1005M Duo"
$ws.Range("D12").Value = "1.9"
$ws.Range("E12").Value = "DDR3-1600 2GB"
$ws.Range("F12").Value = 1021
$ws.Range("G12").Value = 45330

$ws.Range("A13").Value = "Linux"
$ws.Range("B13").Value = "Notebook Lenovo B570"
$ws.Range("C13").Value = "Intel Celeron B800"
$ws.Range("D13").Value = "1.5"
$ws.Range("E13").Value = "DDR3 4GB"
$ws.Range("F13").Value = 1170
$ws.Range("G13").Value = 45323

$ws.Range("A14").Value = "Termux Android "
$ws.Range("B14").Value = "Galaxy Tab A SM-T580"
$ws.Range("C14").Value = "Exynos 7870"
$ws.Range("D14").Value = "1.6"
$ws.Range("E14").Value = "DDR4 2GB"
$ws.Range("F14").Value = 4400
$ws.Range("G14").Value = 45323

$ws.Range("A15").Value = "Termux Android "
$ws.Range("B15").Value = "TV Mi-Box S"
$ws.Range("C15").Value = "Amlogic S905X4-K"
$ws.Range("D15").Value = "1.6"
$ws.Range("E15").Value = "DDR3 2GB"
$ws.Range("F15").Value = 5500
$ws.Range("G15").Value = 45323

$ws.Range("A16").Value = "Termux Android "
$ws.Range("B16").Value = "TV-Box Vontar"
$ws.Range("C16").Value = "Amlogic S905W2"
$ws.Range("D16").Value = "1.8"
$ws.Range("E16").Value = "DDR3L 2GB"
$ws.Range("F16").Value = 15000
$ws.Range("G16").Value = 45323
$ws.Range("H16").Value = "11000-20000"

$ws.Range("A17").Value = "Linux"
$ws.Range("B17").Value = "TV-Box Vontar"
$ws.Range("C17").Value = "Amlogic S905W2"
$ws.Range("D17").Value = "1.8"
$ws.Range("E17").Value = "DDR3 2GB"
$ws.Range("F17").Value = 1700
$ws.Range("G17").Value = 45450
$ws.Range("H17").Value = "1364-1912"

# row 18 - ANTMINER S9 was blank in the benchmark columns; now filled in.
$ws.Range("A18").Value = "Linux"
$ws.Range("B18").Value = "ANTMINER S9"
$ws.Range("C18").Value = "Xilinx ZYNQ XC7Z010"
$ws.Range("D18").Value = "0.766"
$ws.Range("E18").Value = "DDR3L 1GB"
$ws.Range("F18").Value = 8441
$ws.Range("G18").Value = 45465

$ws.Range("A19").Value = "Linux"
$ws.Range("B19").Value = "Libre La Frite"
$ws.Range("C19").Value = "Amlogic S805X-AC"
$ws.Range("D19").Value = "1.2"
$ws.Range("E19").Value = "DDR4-2400 1GB"

$ws.Range("A20").Value = "Linux"
$ws.Range("B20").Value = "Raspberry Pi Zero W"
$ws.Range("C20").Value = "Broadcom BCM2835"
$ws.Range("D20").Value = "1.0"
$ws.Range("E20").Value = "DDR3 512MB"

# row 21 is brand new - reuse G20's (empty, date-formatted) cell style for G21.
$ws.Range("A21").Value = "Linux"
$ws.Range("B21").Value = "Raspberry Pi B"
$ws.Range("C21").Value = "ARM1176JZF-S"
$ws.Range("D21").Value = "0.7"
$ws.Range("E21").Value = "DDR2 256MB"
$ws.Range("G20").Copy($ws.Range("G21"))

# ---- selection follows the new last-edited cell ---------------------------
$ws.Range("F20").Select() | Out-Null
